$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(1,1).Value = 'Parameter'
$ws.Cells.Item(1,2).Value = 'Value'
$ws.Cells.Item(2,1).Value = 'param_TimeStep_starting_index'
$ws.Cells.Item(2,2).Value = 15
$ws.Cells.Item(3,1).Value = 'param_demand1_op_cost_starting_index'
$ws.Cells.Item(3,2).Value = 0
$ws.Cells.Item(4,1).Value = 'param_demand1_inv_cost_starting_index'
$ws.Cells.Item(4,2).Value = 0
$ws.Cells.Item(5,1).Value = 'param_demand2_inv_cost_starting_index'
$ws.Cells.Item(5,2).Value = 0
$ws.Cells.Item(6,1).Value = 'param_demand2_op_cost_starting_index'
$ws.Cells.Item(6,2).Value = 0
$ws.Cells.Item(7,1).Value = 'param_Q_net1_demand2_starting_index'
$ws.Cells.Item(7,2).Value = 1000
$ws.Cells.Item(8,1).Value = 'param_net1_sell_thermal_starting_index'
$ws.Cells.Item(8,2).Value = 0
$ws.Cells.Item(9,1).Value = 'param_net1_buy_electric_starting_index'
$ws.Cells.Item(9,2).Value = 375.8800179973562
$ws.Cells.Item(10,1).Value = 'param_net1_sell_electric_starting_index'
$ws.Cells.Item(10,2).Value = 0
$ws.Cells.Item(11,1).Value = 'param_net1_emissions_starting_index'
$ws.Cells.Item(11,2).Value = 828.1244415259191
$ws.Cells.Item(12,1).Value = 'param_P_net1_bat2_starting_index'
$ws.Cells.Item(12,2).Value = 0
$ws.Cells.Item(13,1).Value = 'param_P_net1_heat_pump2_starting_index'
$ws.Cells.Item(13,2).Value = 19.08749975
$ws.Cells.Item(14,1).Value = 'param_P_net1_bat1_starting_index'
$ws.Cells.Item(14,2).Value = 0
$ws.Cells.Item(15,1).Value = 'param_P_net1_charging_station1_starting_index'
$ws.Cells.Item(15,2).Value = 0
$ws.Cells.Item(16,1).Value = 'param_P_to_net1_starting_index'
$ws.Cells.Item(16,2).Value = 0
$ws.Cells.Item(17,1).Value = 'param_P_net1_demand2_starting_index'
$ws.Cells.Item(17,2).Value = 500
$ws.Cells.Item(18,1).Value = 'param_P_from_net1_starting_index'
$ws.Cells.Item(18,2).Value = 939.7000449933903
$ws.Cells.Item(19,1).Value = 'param_P_net1_demand1_starting_index'
$ws.Cells.Item(19,2).Value = 334.7633794933903
$ws.Cells.Item(20,1).Value = 'param_Q_from_net1_starting_index'
$ws.Cells.Item(20,2).Value = 1257.885068040085
$ws.Cells.Item(21,1).Value = 'param_Q_to_net1_starting_index'
$ws.Cells.Item(21,2).Value = 0
$ws.Cells.Item(22,1).Value = 'param_net1_inv_cost_starting_index'
$ws.Cells.Item(22,2).Value = 0
$ws.Cells.Item(23,1).Value = 'param_P_net1_charging_station2_starting_index'
$ws.Cells.Item(23,2).Value = 66.76166600000001
$ws.Cells.Item(24,1).Value = 'param_P_net1_heat_pump1_starting_index'
$ws.Cells.Item(24,2).Value = 19.08749975
$ws.Cells.Item(25,1).Value = 'param_Q_net1_demand1_starting_index'
$ws.Cells.Item(25,2).Value = 257.8850680400853
$ws.Cells.Item(26,1).Value = 'param_net1_buy_thermal_starting_index'
$ws.Cells.Item(26,2).Value = 465.4174751748316
$ws.Cells.Item(27,1).Value = 'param_net2_buy_electric_starting_index'
$ws.Cells.Item(27,2).Value = 0
$ws.Cells.Item(28,1).Value = 'param_P_net2_bat1_starting_index'
$ws.Cells.Item(28,2).Value = 0
$ws.Cells.Item(29,1).Value = 'param_net2_inv_cost_starting_index'
$ws.Cells.Item(29,2).Value = 0
$ws.Cells.Item(30,1).Value = 'param_net2_buy_thermal_starting_index'
$ws.Cells.Item(30,2).Value = 0
$ws.Cells.Item(31,1).Value = 'param_net2_sell_thermal_starting_index'
$ws.Cells.Item(31,2).Value = 0
$ws.Cells.Item(32,1).Value = 'param_P_net2_bat2_starting_index'
$ws.Cells.Item(32,2).Value = 0
$ws.Cells.Item(33,1).Value = 'param_P_net2_charging_station1_starting_index'
$ws.Cells.Item(33,2).Value = 0
$ws.Cells.Item(34,1).Value = 'param_Q_from_net2_starting_index'
$ws.Cells.Item(34,2).Value = 0
$ws.Cells.Item(35,1).Value = 'param_P_net2_heat_pump2_starting_index'
$ws.Cells.Item(35,2).Value = 0
$ws.Cells.Item(36,1).Value = 'param_P_from_net2_starting_index'
$ws.Cells.Item(36,2).Value = 0
$ws.Cells.Item(37,1).Value = 'param_P_net2_demand2_starting_index'
$ws.Cells.Item(37,2).Value = 0
$ws.Cells.Item(38,1).Value = 'param_Q_net2_demand2_starting_index'
$ws.Cells.Item(38,2).Value = 0
$ws.Cells.Item(39,1).Value = 'param_net2_sell_electric_starting_index'
$ws.Cells.Item(39,2).Value = 0
$ws.Cells.Item(40,1).Value = 'param_P_to_net2_starting_index'
$ws.Cells.Item(40,2).Value = 0
$ws.Cells.Item(41,1).Value = 'param_P_net2_charging_station2_starting_index'
$ws.Cells.Item(41,2).Value = 0
$ws.Cells.Item(42,1).Value = 'param_Q_to_net2_starting_index'
$ws.Cells.Item(42,2).Value = 0
$ws.Cells.Item(43,1).Value = 'param_net2_emissions_starting_index'
$ws.Cells.Item(43,2).Value = 0
$ws.Cells.Item(44,1).Value = 'param_P_net2_heat_pump1_starting_index'
$ws.Cells.Item(44,2).Value = 0
$ws.Cells.Item(45,1).Value = 'param_P_net2_demand1_starting_index'
$ws.Cells.Item(45,2).Value = 0
$ws.Cells.Item(46,1).Value = 'param_Q_net2_demand1_starting_index'
$ws.Cells.Item(46,2).Value = 0
$ws.Cells.Item(47,1).Value = 'param_P_pv1_demand2_starting_index'
$ws.Cells.Item(47,2).Value = 0
$ws.Cells.Item(48,1).Value = 'param_pv1_op_cost_starting_index'
$ws.Cells.Item(48,2).Value = 1
$ws.Cells.Item(49,1).Value = 'param_pv1_emissions_starting_index'
$ws.Cells.Item(49,2).Value = 0.31250025
$ws.Cells.Item(50,1).Value = 'param_P_pv1_demand1_starting_index'
$ws.Cells.Item(50,2).Value = 0
$ws.Cells.Item(51,1).Value = 'param_P_pv1_bat1_starting_index'
$ws.Cells.Item(51,2).Value = 0
$ws.Cells.Item(52,1).Value = 'param_P_pv1_charging_station2_starting_index'
$ws.Cells.Item(52,2).Value = 0
$ws.Cells.Item(53,1).Value = 'param_P_pv1_bat2_starting_index'
$ws.Cells.Item(53,2).Value = 0
$ws.Cells.Item(54,1).Value = 'param_P_pv1_net2_starting_index'
$ws.Cells.Item(54,2).Value = 0
$ws.Cells.Item(55,1).Value = 'param_P_pv1_charging_station1_starting_index'
$ws.Cells.Item(55,2).Value = 0
$ws.Cells.Item(56,1).Value = 'param_P_from_pv1_starting_index'
$ws.Cells.Item(56,2).Value = 0.6250005000000001
$ws.Cells.Item(57,1).Value = 'param_P_pv1_heat_pump2_starting_index'
$ws.Cells.Item(57,2).Value = 0.31250025
$ws.Cells.Item(58,1).Value = 'param_P_pv1_heat_pump1_starting_index'
$ws.Cells.Item(58,2).Value = 0.31250025
$ws.Cells.Item(59,1).Value = 'param_P_pv1_net1_starting_index'
$ws.Cells.Item(59,2).Value = 0
$ws.Cells.Item(60,1).Value = 'param_pv1_inv_cost_starting_index'
$ws.Cells.Item(60,2).Value = 0
$ws.Cells.Item(61,1).Value = 'param_P_pv2_charging_station1_starting_index'
$ws.Cells.Item(61,2).Value = 0
$ws.Cells.Item(62,1).Value = 'param_P_pv2_net1_starting_index'
$ws.Cells.Item(62,2).Value = 0
$ws.Cells.Item(63,1).Value = 'param_P_from_pv2_starting_index'
$ws.Cells.Item(63,2).Value = 1.8
$ws.Cells.Item(64,1).Value = 'param_P_pv2_bat1_starting_index'
$ws.Cells.Item(64,2).Value = 0
$ws.Cells.Item(65,1).Value = 'param_P_pv2_bat2_starting_index'
$ws.Cells.Item(65,2).Value = 0
$ws.Cells.Item(66,1).Value = 'param_P_pv2_demand2_starting_index'
$ws.Cells.Item(66,2).Value = -0
$ws.Cells.Item(67,1).Value = 'param_P_pv2_heat_pump1_starting_index'
$ws.Cells.Item(67,2).Value = 0
$ws.Cells.Item(68,1).Value = 'param_P_pv2_heat_pump2_starting_index'
$ws.Cells.Item(68,2).Value = 0
$ws.Cells.Item(69,1).Value = 'param_pv2_inv_cost_starting_index'
$ws.Cells.Item(69,2).Value = 0
$ws.Cells.Item(70,1).Value = 'param_P_pv2_charging_station2_starting_index'
$ws.Cells.Item(70,2).Value = 1.8
$ws.Cells.Item(71,1).Value = 'param_P_pv2_demand1_starting_index'
$ws.Cells.Item(71,2).Value = 0
$ws.Cells.Item(72,1).Value = 'param_P_pv2_net2_starting_index'
$ws.Cells.Item(72,2).Value = 0
$ws.Cells.Item(73,1).Value = 'param_pv2_op_cost_starting_index'
$ws.Cells.Item(73,2).Value = 1
$ws.Cells.Item(74,1).Value = 'param_pv2_emissions_starting_index'
$ws.Cells.Item(74,2).Value = 0.8999999999999999
$ws.Cells.Item(75,1).Value = 'param_bat1_K_ch_starting_index'
$ws.Cells.Item(75,2).Value = 0
$ws.Cells.Item(76,1).Value = 'param_P_bat1_net1_starting_index'
$ws.Cells.Item(76,2).Value = 0
$ws.Cells.Item(77,1).Value = 'param_P_bat1_heat_pump1_starting_index'
$ws.Cells.Item(77,2).Value = 0
$ws.Cells.Item(78,1).Value = 'param_P_bat1_demand2_starting_index'
$ws.Cells.Item(78,2).Value = 0
$ws.Cells.Item(79,1).Value = 'param_bat1_cumulated_aging_starting_index'
$ws.Cells.Item(79,2).Value = 0
$ws.Cells.Item(80,1).Value = 'param_bat1_op_cost_starting_index'
$ws.Cells.Item(80,2).Value = 1
$ws.Cells.Item(81,1).Value = 'param_bat1_emissions_starting_index'
$ws.Cells.Item(81,2).Value = 0
$ws.Cells.Item(82,1).Value = 'param_bat1_K_dis_starting_index'
$ws.Cells.Item(82,2).Value = 1
$ws.Cells.Item(83,1).Value = 'param_bat1_integer_starting_index'
$ws.Cells.Item(83,2).Value = -0
$ws.Cells.Item(84,1).Value = 'param_bat1_SOC_starting_index'
$ws.Cells.Item(84,2).Value = 0.5
$ws.Cells.Item(85,1).Value = 'param_P_bat1_charging_station1_starting_index'
$ws.Cells.Item(85,2).Value = 0
$ws.Cells.Item(86,1).Value = 'param_P_bat1_heat_pump2_starting_index'
$ws.Cells.Item(86,2).Value = 0
$ws.Cells.Item(87,1).Value = 'param_P_to_bat1_starting_index'
$ws.Cells.Item(87,2).Value = 0
$ws.Cells.Item(88,1).Value = 'param_bat1_inv_cost_starting_index'
$ws.Cells.Item(88,2).Value = 0
$ws.Cells.Item(89,1).Value = 'param_P_from_bat1_starting_index'
$ws.Cells.Item(89,2).Value = 0
$ws.Cells.Item(90,1).Value = 'param_P_bat1_demand1_starting_index'
$ws.Cells.Item(90,2).Value = 0
$ws.Cells.Item(91,1).Value = 'param_P_bat1_charging_station2_starting_index'
$ws.Cells.Item(91,2).Value = 0
$ws.Cells.Item(92,1).Value = 'param_P_bat1_net2_starting_index'
$ws.Cells.Item(92,2).Value = 0
$ws.Cells.Item(93,1).Value = 'param_bat1_SOC_max_starting_index'
$ws.Cells.Item(93,2).Value = 1
$ws.Cells.Item(94,1).Value = 'param_P_bat2_net2_starting_index'
$ws.Cells.Item(94,2).Value = 0
$ws.Cells.Item(95,1).Value = 'param_P_bat2_demand1_starting_index'
$ws.Cells.Item(95,2).Value = 0
$ws.Cells.Item(96,1).Value = 'param_bat2_SOC_starting_index'
$ws.Cells.Item(96,2).Value = 0.5
$ws.Cells.Item(97,1).Value = 'param_P_bat2_charging_station1_starting_index'
$ws.Cells.Item(97,2).Value = 0
$ws.Cells.Item(98,1).Value = 'param_bat2_op_cost_starting_index'
$ws.Cells.Item(98,2).Value = 1
$ws.Cells.Item(99,1).Value = 'param_P_bat2_net1_starting_index'
$ws.Cells.Item(99,2).Value = 0
$ws.Cells.Item(100,1).Value = 'param_bat2_inv_cost_starting_index'
$ws.Cells.Item(100,2).Value = 0
$ws.Cells.Item(101,1).Value = 'param_bat2_emissions_starting_index'
$ws.Cells.Item(101,2).Value = 0
$ws.Cells.Item(102,1).Value = 'param_bat2_cumulated_aging_starting_index'
$ws.Cells.Item(102,2).Value = 0
$ws.Cells.Item(103,1).Value = 'param_P_bat2_demand2_starting_index'
$ws.Cells.Item(103,2).Value = -0
$ws.Cells.Item(104,1).Value = 'param_P_to_bat2_starting_index'
$ws.Cells.Item(104,2).Value = 0
$ws.Cells.Item(105,1).Value = 'param_P_bat2_charging_station2_starting_index'
$ws.Cells.Item(105,2).Value = 0
$ws.Cells.Item(106,1).Value = 'param_P_bat2_heat_pump1_starting_index'
$ws.Cells.Item(106,2).Value = 0
$ws.Cells.Item(107,1).Value = 'param_P_bat2_heat_pump2_starting_index'
$ws.Cells.Item(107,2).Value = 0
$ws.Cells.Item(108,1).Value = 'param_P_from_bat2_starting_index'
$ws.Cells.Item(108,2).Value = 0
$ws.Cells.Item(109,1).Value = 'param_bat2_K_dis_starting_index'
$ws.Cells.Item(109,2).Value = 1
$ws.Cells.Item(110,1).Value = 'param_bat2_K_ch_starting_index'
$ws.Cells.Item(110,2).Value = 0
$ws.Cells.Item(111,1).Value = 'param_bat2_SOC_max_starting_index'
$ws.Cells.Item(111,2).Value = 1
$ws.Cells.Item(112,1).Value = 'param_bat2_integer_starting_index'
$ws.Cells.Item(112,2).Value = -0
$ws.Cells.Item(113,1).Value = 'param_Q_CHP1_demand1_starting_index'
$ws.Cells.Item(113,2).Value = 40
$ws.Cells.Item(114,1).Value = 'param_P_CHP1_bat2_starting_index'
$ws.Cells.Item(114,2).Value = 0
$ws.Cells.Item(115,1).Value = 'param_P_CHP1_charging_station2_starting_index'
$ws.Cells.Item(115,2).Value = 20
$ws.Cells.Item(116,1).Value = 'param_P_CHP1_demand2_starting_index'
$ws.Cells.Item(116,2).Value = (-3.552713678800501 / 1000000000000000)
$ws.Cells.Item(117,1).Value = 'param_CHP1_emissions_starting_index'
$ws.Cells.Item(117,2).Value = 4.83
$ws.Cells.Item(118,1).Value = 'param_P_CHP1_heat_pump1_starting_index'
$ws.Cells.Item(118,2).Value = 0
$ws.Cells.Item(119,1).Value = 'param_P_CHP1_bat1_starting_index'
$ws.Cells.Item(119,2).Value = 0
$ws.Cells.Item(120,1).Value = 'param_P_CHP1_net2_starting_index'
$ws.Cells.Item(120,2).Value = 0
$ws.Cells.Item(121,1).Value = 'param_CHP1_inv_cost_starting_index'
$ws.Cells.Item(121,2).Value = 0
$ws.Cells.Item(122,1).Value = 'param_P_CHP1_net1_starting_index'
$ws.Cells.Item(122,2).Value = 0
$ws.Cells.Item(123,1).Value = 'param_Q_CHP1_net2_starting_index'
$ws.Cells.Item(123,2).Value = 0
$ws.Cells.Item(124,1).Value = 'param_P_CHP1_charging_station1_starting_index'
$ws.Cells.Item(124,2).Value = 0
$ws.Cells.Item(125,1).Value = 'param_P_CHP1_heat_pump2_starting_index'
$ws.Cells.Item(125,2).Value = 0
$ws.Cells.Item(126,1).Value = 'param_P_from_CHP1_starting_index'
$ws.Cells.Item(126,2).Value = 20
$ws.Cells.Item(127,1).Value = 'param_Q_CHP1_net1_starting_index'
$ws.Cells.Item(127,2).Value = 0
$ws.Cells.Item(128,1).Value = 'param_P_CHP1_demand1_starting_index'
$ws.Cells.Item(128,2).Value = 0
$ws.Cells.Item(129,1).Value = 'param_CHP1_fuel_cons_starting_index'
$ws.Cells.Item(129,2).Value = 2.1
$ws.Cells.Item(130,1).Value = 'param_CHP1_op_cost_starting_index'
$ws.Cells.Item(130,2).Value = 10.5
$ws.Cells.Item(131,1).Value = 'param_Q_from_CHP1_starting_index'
$ws.Cells.Item(131,2).Value = 40
$ws.Cells.Item(132,1).Value = 'param_Q_CHP1_demand2_starting_index'
$ws.Cells.Item(132,2).Value = 0
$ws.Cells.Item(133,1).Value = 'param_P_CHP2_net2_starting_index'
$ws.Cells.Item(133,2).Value = 0
$ws.Cells.Item(134,1).Value = 'param_P_CHP2_bat2_starting_index'
$ws.Cells.Item(134,2).Value = 0
$ws.Cells.Item(135,1).Value = 'param_P_from_CHP2_starting_index'
$ws.Cells.Item(135,2).Value = 20
$ws.Cells.Item(136,1).Value = 'param_CHP2_op_cost_starting_index'
$ws.Cells.Item(136,2).Value = 10.5
$ws.Cells.Item(137,1).Value = 'param_Q_CHP2_net1_starting_index'
$ws.Cells.Item(137,2).Value = 0
$ws.Cells.Item(138,1).Value = 'param_P_CHP2_heat_pump2_starting_index'
$ws.Cells.Item(138,2).Value = 0
$ws.Cells.Item(139,1).Value = 'param_Q_CHP2_net2_starting_index'
$ws.Cells.Item(139,2).Value = 0
$ws.Cells.Item(140,1).Value = 'param_Q_CHP2_demand1_starting_index'
$ws.Cells.Item(140,2).Value = 40
$ws.Cells.Item(141,1).Value = 'param_CHP2_fuel_cons_starting_index'
$ws.Cells.Item(141,2).Value = 2.1
$ws.Cells.Item(142,1).Value = 'param_P_CHP2_charging_station1_starting_index'
$ws.Cells.Item(142,2).Value = 0
$ws.Cells.Item(143,1).Value = 'param_CHP2_inv_cost_starting_index'
$ws.Cells.Item(143,2).Value = 0
$ws.Cells.Item(144,1).Value = 'param_P_CHP2_bat1_starting_index'
$ws.Cells.Item(144,2).Value = 0
$ws.Cells.Item(145,1).Value = 'param_P_CHP2_demand1_starting_index'
$ws.Cells.Item(145,2).Value = 0
$ws.Cells.Item(146,1).Value = 'param_Q_from_CHP2_starting_index'
$ws.Cells.Item(146,2).Value = 40
$ws.Cells.Item(147,1).Value = 'param_Q_CHP2_demand2_starting_index'
$ws.Cells.Item(147,2).Value = 0
$ws.Cells.Item(148,1).Value = 'param_P_CHP2_heat_pump1_starting_index'
$ws.Cells.Item(148,2).Value = 0
$ws.Cells.Item(149,1).Value = 'param_P_CHP2_net1_starting_index'
$ws.Cells.Item(149,2).Value = 0
$ws.Cells.Item(150,1).Value = 'param_CHP2_emissions_starting_index'
$ws.Cells.Item(150,2).Value = 4.83
$ws.Cells.Item(151,1).Value = 'param_P_CHP2_demand2_starting_index'
$ws.Cells.Item(151,2).Value = (-3.552713678800501 / 1000000000000000)
$ws.Cells.Item(152,1).Value = 'param_P_CHP2_charging_station2_starting_index'
$ws.Cells.Item(152,2).Value = 20
$ws.Cells.Item(153,1).Value = 'param_Q_solar_th1_net2_starting_index'
$ws.Cells.Item(153,2).Value = 0
$ws.Cells.Item(154,1).Value = 'param_Q_from_solar_th1_starting_index'
$ws.Cells.Item(154,2).Value = 0.416667
$ws.Cells.Item(155,1).Value = 'param_Q_solar_th1_demand1_starting_index'
$ws.Cells.Item(155,2).Value = 0.416667
$ws.Cells.Item(156,1).Value = 'param_Q_solar_th1_net1_starting_index'
$ws.Cells.Item(156,2).Value = 0
$ws.Cells.Item(157,1).Value = 'param_solar_th1_op_cost_starting_index'
$ws.Cells.Item(157,2).Value = 1
$ws.Cells.Item(158,1).Value = 'param_solar_th1_inv_cost_starting_index'
$ws.Cells.Item(158,2).Value = 0
$ws.Cells.Item(159,1).Value = 'param_solar_th1_emissions_starting_index'
$ws.Cells.Item(159,2).Value = 0.2083335
$ws.Cells.Item(160,1).Value = 'param_Q_solar_th1_demand2_starting_index'
$ws.Cells.Item(160,2).Value = 0
$ws.Cells.Item(161,1).Value = 'param_Q_solar_th2_net1_starting_index'
$ws.Cells.Item(161,2).Value = 0
$ws.Cells.Item(162,1).Value = 'param_Q_solar_th2_net2_starting_index'
$ws.Cells.Item(162,2).Value = 0
$ws.Cells.Item(163,1).Value = 'param_Q_solar_th2_demand2_starting_index'
$ws.Cells.Item(163,2).Value = 0
$ws.Cells.Item(164,1).Value = 'param_solar_th2_emissions_starting_index'
$ws.Cells.Item(164,2).Value = 0.6000000000000001
$ws.Cells.Item(165,1).Value = 'param_Q_from_solar_th2_starting_index'
$ws.Cells.Item(165,2).Value = 1.2
$ws.Cells.Item(166,1).Value = 'param_solar_th2_op_cost_starting_index'
$ws.Cells.Item(166,2).Value = 1
$ws.Cells.Item(167,1).Value = 'param_Q_solar_th2_demand1_starting_index'
$ws.Cells.Item(167,2).Value = 1.2
$ws.Cells.Item(168,1).Value = 'param_solar_th2_inv_cost_starting_index'
$ws.Cells.Item(168,2).Value = 0
$ws.Cells.Item(169,1).Value = 'param_P_pvt1_bat2_starting_index'
$ws.Cells.Item(169,2).Value = 0
$ws.Cells.Item(170,1).Value = 'param_P_from_pvt1_starting_index'
$ws.Cells.Item(170,2).Value = 0.833334
$ws.Cells.Item(171,1).Value = 'param_pvt1_emissions_starting_index'
$ws.Cells.Item(171,2).Value = 0.5416671000000001
$ws.Cells.Item(172,1).Value = 'param_pvt1_inv_cost_starting_index'
$ws.Cells.Item(172,2).Value = 0
$ws.Cells.Item(173,1).Value = 'param_P_pvt1_bat1_starting_index'
$ws.Cells.Item(173,2).Value = 0
$ws.Cells.Item(174,1).Value = 'param_P_pvt1_net2_starting_index'
$ws.Cells.Item(174,2).Value = 0
$ws.Cells.Item(175,1).Value = 'param_pvt1_op_cost_starting_index'
$ws.Cells.Item(175,2).Value = 1
$ws.Cells.Item(176,1).Value = 'param_P_pvt1_heat_pump2_starting_index'
$ws.Cells.Item(176,2).Value = 0
$ws.Cells.Item(177,1).Value = 'param_P_pvt1_charging_station2_starting_index'
$ws.Cells.Item(177,2).Value = 0.833334
$ws.Cells.Item(178,1).Value = 'param_Q_pvt1_net1_starting_index'
$ws.Cells.Item(178,2).Value = 0
$ws.Cells.Item(179,1).Value = 'param_Q_from_pvt1_starting_index'
$ws.Cells.Item(179,2).Value = 1.0833342
$ws.Cells.Item(180,1).Value = 'param_P_pvt1_net1_starting_index'
$ws.Cells.Item(180,2).Value = 0
$ws.Cells.Item(181,1).Value = 'param_Q_pvt1_demand2_starting_index'
$ws.Cells.Item(181,2).Value = 0
$ws.Cells.Item(182,1).Value = 'param_P_pvt1_demand1_starting_index'
$ws.Cells.Item(182,2).Value = 0
$ws.Cells.Item(183,1).Value = 'param_P_pvt1_heat_pump1_starting_index'
$ws.Cells.Item(183,2).Value = 0
$ws.Cells.Item(184,1).Value = 'param_Q_pvt1_net2_starting_index'
$ws.Cells.Item(184,2).Value = 0
$ws.Cells.Item(185,1).Value = 'param_P_pvt1_charging_station1_starting_index'
$ws.Cells.Item(185,2).Value = 0
$ws.Cells.Item(186,1).Value = 'param_Q_pvt1_demand1_starting_index'
$ws.Cells.Item(186,2).Value = 1.0833342
$ws.Cells.Item(187,1).Value = 'param_P_pvt1_demand2_starting_index'
$ws.Cells.Item(187,2).Value = 0
$ws.Cells.Item(188,1).Value = 'param_P_pvt2_bat1_starting_index'
$ws.Cells.Item(188,2).Value = 0
$ws.Cells.Item(189,1).Value = 'param_P_pvt2_demand2_starting_index'
$ws.Cells.Item(189,2).Value = 0
$ws.Cells.Item(190,1).Value = 'param_P_pvt2_net1_starting_index'
$ws.Cells.Item(190,2).Value = 0
$ws.Cells.Item(191,1).Value = 'param_Q_from_pvt2_starting_index'
$ws.Cells.Item(191,2).Value = 1.56
$ws.Cells.Item(192,1).Value = 'param_P_from_pvt2_starting_index'
$ws.Cells.Item(192,2).Value = 1.2
$ws.Cells.Item(193,1).Value = 'param_P_pvt2_charging_station2_starting_index'
$ws.Cells.Item(193,2).Value = 0
$ws.Cells.Item(194,1).Value = 'param_P_pvt2_heat_pump2_starting_index'
$ws.Cells.Item(194,2).Value = 0.6000000000000001
$ws.Cells.Item(195,1).Value = 'param_P_pvt2_bat2_starting_index'
$ws.Cells.Item(195,2).Value = 0
$ws.Cells.Item(196,1).Value = 'param_pvt2_emissions_starting_index'
$ws.Cells.Item(196,2).Value = 0.7800000000000001
$ws.Cells.Item(197,1).Value = 'param_Q_pvt2_demand1_starting_index'
$ws.Cells.Item(197,2).Value = 1.56
$ws.Cells.Item(198,1).Value = 'param_P_pvt2_charging_station1_starting_index'
$ws.Cells.Item(198,2).Value = 0
$ws.Cells.Item(199,1).Value = 'param_P_pvt2_demand1_starting_index'
$ws.Cells.Item(199,2).Value = 0
$ws.Cells.Item(200,1).Value = 'param_Q_pvt2_net2_starting_index'
$ws.Cells.Item(200,2).Value = 0
$ws.Cells.Item(201,1).Value = 'param_P_pvt2_net2_starting_index'
$ws.Cells.Item(201,2).Value = 0
$ws.Cells.Item(202,1).Value = 'param_P_pvt2_heat_pump1_starting_index'
$ws.Cells.Item(202,2).Value = 0.6000000000000001
$ws.Cells.Item(203,1).Value = 'param_Q_pvt2_net1_starting_index'
$ws.Cells.Item(203,2).Value = 0
$ws.Cells.Item(204,1).Value = 'param_pvt2_inv_cost_starting_index'
$ws.Cells.Item(204,2).Value = 0
$ws.Cells.Item(205,1).Value = 'param_Q_pvt2_demand2_starting_index'
$ws.Cells.Item(205,2).Value = 0
$ws.Cells.Item(206,1).Value = 'param_pvt2_op_cost_starting_index'
$ws.Cells.Item(206,2).Value = 1
$ws.Cells.Item(207,1).Value = 'param_charging_station1_inv_cost_starting_index'
$ws.Cells.Item(207,2).Value = 0
$ws.Cells.Item(208,1).Value = 'param_charging_station1_op_cost_starting_index'
$ws.Cells.Item(208,2).Value = 0
$ws.Cells.Item(209,1).Value = 'param_charging_station1_emissions_starting_index'
$ws.Cells.Item(209,2).Value = 0
$ws.Cells.Item(210,1).Value = 'param_charging_station2_emissions_starting_index'
$ws.Cells.Item(210,2).Value = 5.46975
$ws.Cells.Item(211,1).Value = 'param_charging_station2_inv_cost_starting_index'
$ws.Cells.Item(211,2).Value = 0
$ws.Cells.Item(212,1).Value = 'param_charging_station2_op_cost_starting_index'
$ws.Cells.Item(212,2).Value = -65.637
$ws.Cells.Item(213,1).Value = 'param_Q_heat_pump1_net1_starting_index'
$ws.Cells.Item(213,2).Value = 0
$ws.Cells.Item(214,1).Value = 'param_heat_pump1_op_cost_starting_index'
$ws.Cells.Item(214,2).Value = 8.561643835616438
$ws.Cells.Item(215,1).Value = 'param_heat_pump1_emissions_starting_index'
$ws.Cells.Item(215,2).Value = 2.76
$ws.Cells.Item(216,1).Value = 'param_Q_from_heat_pump1_starting_index'
$ws.Cells.Item(216,2).Value = 80
$ws.Cells.Item(217,1).Value = 'param_Q_to_heat_pump1_starting_index'
$ws.Cells.Item(217,2).Value = 0
$ws.Cells.Item(218,1).Value = 'param_P_from_heat_pump1_starting_index'
$ws.Cells.Item(218,2).Value = 0
$ws.Cells.Item(219,1).Value = 'param_heat_pump1_inv_cost_starting_index'
$ws.Cells.Item(219,2).Value = 0
$ws.Cells.Item(220,1).Value = 'param_Q_heat_pump1_net2_starting_index'
$ws.Cells.Item(220,2).Value = 0
$ws.Cells.Item(221,1).Value = 'param_Q_heat_pump1_demand1_starting_index'
$ws.Cells.Item(221,2).Value = 80
$ws.Cells.Item(222,1).Value = 'param_Q_heat_pump1_demand2_starting_index'
$ws.Cells.Item(222,2).Value = 0
$ws.Cells.Item(223,1).Value = 'param_P_to_heat_pump1_starting_index'
$ws.Cells.Item(223,2).Value = 20
$ws.Cells.Item(224,1).Value = 'param_P_to_heat_pump2_starting_index'
$ws.Cells.Item(224,2).Value = 20
$ws.Cells.Item(225,1).Value = 'param_Q_from_heat_pump2_starting_index'
$ws.Cells.Item(225,2).Value = 80
$ws.Cells.Item(226,1).Value = 'param_heat_pump2_emissions_starting_index'
$ws.Cells.Item(226,2).Value = 2.76
$ws.Cells.Item(227,1).Value = 'param_P_from_heat_pump2_starting_index'
$ws.Cells.Item(227,2).Value = 0
$ws.Cells.Item(228,1).Value = 'param_Q_heat_pump2_demand2_starting_index'
$ws.Cells.Item(228,2).Value = 0
$ws.Cells.Item(229,1).Value = 'param_Q_heat_pump2_net1_starting_index'
$ws.Cells.Item(229,2).Value = 0
$ws.Cells.Item(230,1).Value = 'param_Q_heat_pump2_net2_starting_index'
$ws.Cells.Item(230,2).Value = 0
$ws.Cells.Item(231,1).Value = 'param_Q_heat_pump2_demand1_starting_index'
$ws.Cells.Item(231,2).Value = 80
$ws.Cells.Item(232,1).Value = 'param_heat_pump2_op_cost_starting_index'
$ws.Cells.Item(232,2).Value = 8.561643835616438
$ws.Cells.Item(233,1).Value = 'param_heat_pump2_inv_cost_starting_index'
$ws.Cells.Item(233,2).Value = 0
$ws.Cells.Item(234,1).Value = 'param_Q_to_heat_pump2_starting_index'
$ws.Cells.Item(234,2).Value = 0
$ws.Cells.Item(235,1).Value = 'param_total_emissions_starting_index'
$ws.Cells.Item(235,2).Value = 849.9866917759191
$ws.Cells.Item(236,1).Value = 'param_total_sell_starting_index'
$ws.Cells.Item(236,2).Value = 0
$ws.Cells.Item(237,1).Value = 'param_total_buy_starting_index'
$ws.Cells.Item(237,2).Value = 841.2974931721877
$ws.Cells.Item(238,1).Value = 'param_total_operation_cost_starting_index'
$ws.Cells.Item(238,2).Value = -23.51371232876713
